# Taetigkeitsdokumentation.xlsx - "Add files via upload"
#
# Adds two new activity-log rows (21 and 22) on the "Tabelle1" sheet for
# the game_launcher.py functions `open_regelwerk` and `reset_regelwerk`,
# each attributed to Niklas. The dependent SUM/ratio formulas elsewhere
# on the sheet (L2, H3:M3, ...) recalculate automatically once the new
# source values are in place.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Row 21: open_regelwerk, 4 lines, weight 1, author Niklas
$ws.Range("A21").Value = "game_launcher.py"
$ws.Range("B21").Value = "open_regelwerk"
$ws.Range("C21").Value = 4
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "Niklas"

# Row 22: reset_regelwerk, 2 lines, weight 1, author Niklas
$ws.Range("A22").Value = "game_launcher.py"
$ws.Range("B22").Value = "reset_regelwerk"
$ws.Range("C22").Value = 2
$ws.Range("D22").Value = 1
$ws.Range("E22").Value = "Niklas"

# Recalculate so every dependent formula carries a fresh cached value.
$excel.CalculateFull()

# Restore the view to match the saved workbook state: scrolled back to
# the top and with H17 as the active selection (instead of D23).
$ws.Range("H17").Select() | Out-Null
